$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("D3").Value = "Circuitos Elétricos 2"
$ws.Range("C4").Value = "Circuitos Elétricos 2"
$ws.Range("C6").Value = "Circuitos Elétricos 2"
$ws.Range("E6").Value = "-"
$ws.Range("C7").Value = "-"
